$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '50.989.60'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '3.042.98'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '388.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.529'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.78%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.577'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.54'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0843'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.03%  '
$ws.Range("D13").Value = '3.530.74'
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("D16").Value = '3.047.45'
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.992'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").Value = '51.035.62'
$ws.Range("E19").Value = '  -1.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.90%  '
$ws.Range("D22").Value = '0.0₃0950'
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.44'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '262.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.84'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.13'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.161'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.104'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.45'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0483'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '35.68'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.89%  '
$ws.Range("B35").Value = 'Toncoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '49.93'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.11%  '
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.30'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.287'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '128.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.45'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.53%  '
$ws.Range("E42").Value = '  -2.68%  '
$ws.Range("E43").Value = '  -1.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.71'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.86%  '
$ws.Range("E48").Value = '  -0.71%  '
$ws.Range("D49").Value = '2.057.29'
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("B50").Value = 'FlareNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/2hOSU_JYX+flarenetwork-flr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0474'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +20.74%  '
$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0313'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.99%  '
